$wb = $excel.ActiveWorkbook

# --- PIR sheet: append rows 281-293 ---
$wsPIR = $wb.Worksheets.Item("PIR")
$wsPIR.Cells.Item(281, 1).NumberFormat = "@"
$wsPIR.Cells.Item(281, 1).Value = "2026-01-28"
$wsPIR.Cells.Item(281, 2).Value = "12:27:54"
$wsPIR.Cells.Item(281, 3).Value = "12:00"
$wsPIR.Cells.Item(281, 4).Value = "Bathroom"
$wsPIR.Cells.Item(281, 5).Value = "No Motion"
$wsPIR.Cells.Item(281, 6).Value = "Inactive"
$wsPIR.Cells.Item(282, 1).NumberFormat = "@"
$wsPIR.Cells.Item(282, 1).Value = "2026-01-28"
$wsPIR.Cells.Item(282, 2).Value = "12:27:57"
$wsPIR.Cells.Item(282, 3).Value = "12:00"
$wsPIR.Cells.Item(282, 4).Value = "Bathroom"
$wsPIR.Cells.Item(282, 5).Value = "No Motion"
$wsPIR.Cells.Item(282, 6).Value = "Inactive"
$wsPIR.Cells.Item(283, 1).NumberFormat = "@"
$wsPIR.Cells.Item(283, 1).Value = "2026-01-28"
$wsPIR.Cells.Item(283, 2).Value = "12:28:00"
$wsPIR.Cells.Item(283, 3).Value = "12:00"
$wsPIR.Cells.Item(283, 4).Value = "Bathroom"
$wsPIR.Cells.Item(283, 5).Value = "No Motion"
$wsPIR.Cells.Item(283, 6).Value = "Inactive"
$wsPIR.Cells.Item(284, 1).NumberFormat = "@"
$wsPIR.Cells.Item(284, 1).Value = "2026-01-28"
$wsPIR.Cells.Item(284, 2).Value = "12:28:05"
$wsPIR.Cells.Item(284, 3).Value = "12:00"
$wsPIR.Cells.Item(284, 4).Value = "Bathroom"
$wsPIR.Cells.Item(284, 5).Value = "No Motion"
$wsPIR.Cells.Item(284, 6).Value = "Inactive"
$wsPIR.Cells.Item(285, 1).NumberFormat = "@"
$wsPIR.Cells.Item(285, 1).Value = "2026-01-28"
$wsPIR.Cells.Item(285, 2).Value = "12:28:10"
$wsPIR.Cells.Item(285, 3).Value = "12:00"
$wsPIR.Cells.Item(285, 4).Value = "Bathroom"
$wsPIR.Cells.Item(285, 5).Value = "No Motion"
$wsPIR.Cells.Item(285, 6).Value = "Inactive"
$wsPIR.Cells.Item(286, 1).NumberFormat = "@"
$wsPIR.Cells.Item(286, 1).Value = "2026-01-28"
$wsPIR.Cells.Item(286, 2).Value = "12:28:17"
$wsPIR.Cells.Item(286, 3).Value = "12:00"
$wsPIR.Cells.Item(286, 4).Value = "Bathroom"
$wsPIR.Cells.Item(286, 5).Value = "No Motion"
$wsPIR.Cells.Item(286, 6).Value = "Inactive"
$wsPIR.Cells.Item(287, 1).NumberFormat = "@"
$wsPIR.Cells.Item(287, 1).Value = "2026-01-28"
$wsPIR.Cells.Item(287, 2).Value = "12:28:21"
$wsPIR.Cells.Item(287, 3).Value = "12:00"
$wsPIR.Cells.Item(287, 4).Value = "Bathroom"
$wsPIR.Cells.Item(287, 5).Value = "No Motion"
$wsPIR.Cells.Item(287, 6).Value = "Inactive"
$wsPIR.Cells.Item(288, 1).NumberFormat = "@"
$wsPIR.Cells.Item(288, 1).Value = "2026-01-28"
$wsPIR.Cells.Item(288, 2).Value = "12:28:25"
$wsPIR.Cells.Item(288, 3).Value = "12:00"
$wsPIR.Cells.Item(288, 4).Value = "Bathroom"
$wsPIR.Cells.Item(288, 5).Value = "No Motion"
$wsPIR.Cells.Item(288, 6).Value = "Inactive"
$wsPIR.Cells.Item(289, 1).NumberFormat = "@"
$wsPIR.Cells.Item(289, 1).Value = "2026-01-28"
$wsPIR.Cells.Item(289, 2).Value = "12:28:30"
$wsPIR.Cells.Item(289, 3).Value = "12:00"
$wsPIR.Cells.Item(289, 4).Value = "Bathroom"
$wsPIR.Cells.Item(289, 5).Value = "No Motion"
$wsPIR.Cells.Item(289, 6).Value = "Inactive"
$wsPIR.Cells.Item(290, 1).NumberFormat = "@"
$wsPIR.Cells.Item(290, 1).Value = "2026-01-28"
$wsPIR.Cells.Item(290, 2).Value = "12:28:37"
$wsPIR.Cells.Item(290, 3).Value = "12:00"
$wsPIR.Cells.Item(290, 4).Value = "Bathroom"
$wsPIR.Cells.Item(290, 5).Value = "No Motion"
$wsPIR.Cells.Item(290, 6).Value = "Inactive"
$wsPIR.Cells.Item(291, 1).NumberFormat = "@"
$wsPIR.Cells.Item(291, 1).Value = "2026-01-28"
$wsPIR.Cells.Item(291, 2).Value = "12:28:41"
$wsPIR.Cells.Item(291, 3).Value = "12:00"
$wsPIR.Cells.Item(291, 4).Value = "Bathroom"
$wsPIR.Cells.Item(291, 5).Value = "No Motion"
$wsPIR.Cells.Item(291, 6).Value = "Inactive"
$wsPIR.Cells.Item(292, 1).NumberFormat = "@"
$wsPIR.Cells.Item(292, 1).Value = "2026-01-28"
$wsPIR.Cells.Item(292, 2).Value = "12:28:45"
$wsPIR.Cells.Item(292, 3).Value = "12:00"
$wsPIR.Cells.Item(292, 4).Value = "Bathroom"
$wsPIR.Cells.Item(292, 5).Value = "No Motion"
$wsPIR.Cells.Item(292, 6).Value = "Inactive"
$wsPIR.Cells.Item(293, 1).NumberFormat = "@"
$wsPIR.Cells.Item(293, 1).Value = "2026-01-28"
$wsPIR.Cells.Item(293, 2).Value = "12:28:50"
$wsPIR.Cells.Item(293, 3).Value = "12:00"
$wsPIR.Cells.Item(293, 4).Value = "Bathroom"
$wsPIR.Cells.Item(293, 5).Value = "No Motion"
$wsPIR.Cells.Item(293, 6).Value = "Inactive"

# --- Humidity sheet: append rows 261-274 ---
$wsHum = $wb.Worksheets.Item("Humidity")
$wsHum.Cells.Item(261, 1).NumberFormat = "@"
$wsHum.Cells.Item(261, 1).Value = "2026-01-28"
$wsHum.Cells.Item(261, 2).Value = "12:27:54"
$wsHum.Cells.Item(261, 3).Value = "12:00"
$wsHum.Cells.Item(261, 4).Value = "Bathroom"
$wsHum.Cells.Item(261, 5).NumberFormat = "@"
$wsHum.Cells.Item(261, 5).Value = "86.6%"
$wsHum.Cells.Item(261, 6).Value = "Active"
$wsHum.Cells.Item(262, 1).NumberFormat = "@"
$wsHum.Cells.Item(262, 1).Value = "2026-01-28"
$wsHum.Cells.Item(262, 2).Value = "12:27:56"
$wsHum.Cells.Item(262, 3).Value = "12:00"
$wsHum.Cells.Item(262, 4).Value = "Bathroom"
$wsHum.Cells.Item(262, 5).NumberFormat = "@"
$wsHum.Cells.Item(262, 5).Value = "87.6%"
$wsHum.Cells.Item(262, 6).Value = "Active"
$wsHum.Cells.Item(263, 1).NumberFormat = "@"
$wsHum.Cells.Item(263, 1).Value = "2026-01-28"
$wsHum.Cells.Item(263, 2).Value = "12:27:59"
$wsHum.Cells.Item(263, 3).Value = "12:00"
$wsHum.Cells.Item(263, 4).Value = "Bathroom"
$wsHum.Cells.Item(263, 5).NumberFormat = "@"
$wsHum.Cells.Item(263, 5).Value = "87.6%"
$wsHum.Cells.Item(263, 6).Value = "Active"
$wsHum.Cells.Item(264, 1).NumberFormat = "@"
$wsHum.Cells.Item(264, 1).Value = "2026-01-28"
$wsHum.Cells.Item(264, 2).Value = "12:28:03"
$wsHum.Cells.Item(264, 3).Value = "12:00"
$wsHum.Cells.Item(264, 4).Value = "Bathroom"
$wsHum.Cells.Item(264, 5).NumberFormat = "@"
$wsHum.Cells.Item(264, 5).Value = "86.7%"
$wsHum.Cells.Item(264, 6).Value = "Active"
$wsHum.Cells.Item(265, 1).NumberFormat = "@"
$wsHum.Cells.Item(265, 1).Value = "2026-01-28"
$wsHum.Cells.Item(265, 2).Value = "12:28:07"
$wsHum.Cells.Item(265, 3).Value = "12:00"
$wsHum.Cells.Item(265, 4).Value = "Bathroom"
$wsHum.Cells.Item(265, 5).NumberFormat = "@"
$wsHum.Cells.Item(265, 5).Value = "87.6%"
$wsHum.Cells.Item(265, 6).Value = "Active"
$wsHum.Cells.Item(266, 1).NumberFormat = "@"
$wsHum.Cells.Item(266, 1).Value = "2026-01-28"
$wsHum.Cells.Item(266, 2).Value = "12:28:15"
$wsHum.Cells.Item(266, 3).Value = "12:00"
$wsHum.Cells.Item(266, 4).Value = "Bathroom"
$wsHum.Cells.Item(266, 5).NumberFormat = "@"
$wsHum.Cells.Item(266, 5).Value = "87.6%"
$wsHum.Cells.Item(266, 6).Value = "Active"
$wsHum.Cells.Item(267, 1).NumberFormat = "@"
$wsHum.Cells.Item(267, 1).Value = "2026-01-28"
$wsHum.Cells.Item(267, 2).Value = "12:28:19"
$wsHum.Cells.Item(267, 3).Value = "12:00"
$wsHum.Cells.Item(267, 4).Value = "Bathroom"
$wsHum.Cells.Item(267, 5).NumberFormat = "@"
$wsHum.Cells.Item(267, 5).Value = "87.6%"
$wsHum.Cells.Item(267, 6).Value = "Active"
$wsHum.Cells.Item(268, 1).NumberFormat = "@"
$wsHum.Cells.Item(268, 1).Value = "2026-01-28"
$wsHum.Cells.Item(268, 2).Value = "12:28:23"
$wsHum.Cells.Item(268, 3).Value = "12:00"
$wsHum.Cells.Item(268, 4).Value = "Bathroom"
$wsHum.Cells.Item(268, 5).NumberFormat = "@"
$wsHum.Cells.Item(268, 5).Value = "86.7%"
$wsHum.Cells.Item(268, 6).Value = "Active"
$wsHum.Cells.Item(269, 1).NumberFormat = "@"
$wsHum.Cells.Item(269, 1).Value = "2026-01-28"
$wsHum.Cells.Item(269, 2).Value = "12:28:27"
$wsHum.Cells.Item(269, 3).Value = "12:00"
$wsHum.Cells.Item(269, 4).Value = "Bathroom"
$wsHum.Cells.Item(269, 5).NumberFormat = "@"
$wsHum.Cells.Item(269, 5).Value = "87.6%"
$wsHum.Cells.Item(269, 6).Value = "Active"
$wsHum.Cells.Item(270, 1).NumberFormat = "@"
$wsHum.Cells.Item(270, 1).Value = "2026-01-28"
$wsHum.Cells.Item(270, 2).Value = "12:28:35"
$wsHum.Cells.Item(270, 3).Value = "12:00"
$wsHum.Cells.Item(270, 4).Value = "Bathroom"
$wsHum.Cells.Item(270, 5).NumberFormat = "@"
$wsHum.Cells.Item(270, 5).Value = "87.6%"
$wsHum.Cells.Item(270, 6).Value = "Active"
$wsHum.Cells.Item(271, 1).NumberFormat = "@"
$wsHum.Cells.Item(271, 1).Value = "2026-01-28"
$wsHum.Cells.Item(271, 2).Value = "12:28:39"
$wsHum.Cells.Item(271, 3).Value = "12:00"
$wsHum.Cells.Item(271, 4).Value = "Bathroom"
$wsHum.Cells.Item(271, 5).NumberFormat = "@"
$wsHum.Cells.Item(271, 5).Value = "87.6%"
$wsHum.Cells.Item(271, 6).Value = "Active"
$wsHum.Cells.Item(272, 1).NumberFormat = "@"
$wsHum.Cells.Item(272, 1).Value = "2026-01-28"
$wsHum.Cells.Item(272, 2).Value = "12:28:43"
$wsHum.Cells.Item(272, 3).Value = "12:00"
$wsHum.Cells.Item(272, 4).Value = "Bathroom"
$wsHum.Cells.Item(272, 5).NumberFormat = "@"
$wsHum.Cells.Item(272, 5).Value = "86.7%"
$wsHum.Cells.Item(272, 6).Value = "Active"
$wsHum.Cells.Item(273, 1).NumberFormat = "@"
$wsHum.Cells.Item(273, 1).Value = "2026-01-28"
$wsHum.Cells.Item(273, 2).Value = "12:28:47"
$wsHum.Cells.Item(273, 3).Value = "12:00"
$wsHum.Cells.Item(273, 4).Value = "Bathroom"
$wsHum.Cells.Item(273, 5).NumberFormat = "@"
$wsHum.Cells.Item(273, 5).Value = "87.6%"
$wsHum.Cells.Item(273, 6).Value = "Active"
$wsHum.Cells.Item(274, 1).NumberFormat = "@"
$wsHum.Cells.Item(274, 1).Value = "2026-01-28"
$wsHum.Cells.Item(274, 2).Value = "12:28:51"
$wsHum.Cells.Item(274, 3).Value = "12:00"
$wsHum.Cells.Item(274, 4).Value = "Bathroom"
$wsHum.Cells.Item(274, 5).NumberFormat = "@"
$wsHum.Cells.Item(274, 5).Value = "86.7%"
$wsHum.Cells.Item(274, 6).Value = "Active"

# --- Temperature sheet: append rows 261-274 ---
$wsTemp = $wb.Worksheets.Item("Temperature")
$wsTemp.Cells.Item(261, 1).NumberFormat = "@"
$wsTemp.Cells.Item(261, 1).Value = "2026-01-28"
$wsTemp.Cells.Item(261, 2).Value = "12:27:55"
$wsTemp.Cells.Item(261, 3).Value = "12:00"
$wsTemp.Cells.Item(261, 4).Value = "Bathroom"
$wsTemp.Cells.Item(261, 5).Value = "22.9C"
$wsTemp.Cells.Item(261, 6).Value = "Active"
$wsTemp.Cells.Item(262, 1).NumberFormat = "@"
$wsTemp.Cells.Item(262, 1).Value = "2026-01-28"
$wsTemp.Cells.Item(262, 2).Value = "12:27:57"
$wsTemp.Cells.Item(262, 3).Value = "12:00"
$wsTemp.Cells.Item(262, 4).Value = "Bathroom"
$wsTemp.Cells.Item(262, 5).Value = "22.9C"
$wsTemp.Cells.Item(262, 6).Value = "Active"
$wsTemp.Cells.Item(263, 1).NumberFormat = "@"
$wsTemp.Cells.Item(263, 1).Value = "2026-01-28"
$wsTemp.Cells.Item(263, 2).Value = "12:28:00"
$wsTemp.Cells.Item(263, 3).Value = "12:00"
$wsTemp.Cells.Item(263, 4).Value = "Bathroom"
$wsTemp.Cells.Item(263, 5).Value = "22.9C"
$wsTemp.Cells.Item(263, 6).Value = "Active"
$wsTemp.Cells.Item(264, 1).NumberFormat = "@"
$wsTemp.Cells.Item(264, 1).Value = "2026-01-28"
$wsTemp.Cells.Item(264, 2).Value = "12:28:04"
$wsTemp.Cells.Item(264, 3).Value = "12:00"
$wsTemp.Cells.Item(264, 4).Value = "Bathroom"
$wsTemp.Cells.Item(264, 5).Value = "22.9C"
$wsTemp.Cells.Item(264, 6).Value = "Active"
$wsTemp.Cells.Item(265, 1).NumberFormat = "@"
$wsTemp.Cells.Item(265, 1).Value = "2026-01-28"
$wsTemp.Cells.Item(265, 2).Value = "12:28:08"
$wsTemp.Cells.Item(265, 3).Value = "12:00"
$wsTemp.Cells.Item(265, 4).Value = "Bathroom"
$wsTemp.Cells.Item(265, 5).Value = "22.9C"
$wsTemp.Cells.Item(265, 6).Value = "Active"
$wsTemp.Cells.Item(266, 1).NumberFormat = "@"
$wsTemp.Cells.Item(266, 1).Value = "2026-01-28"
$wsTemp.Cells.Item(266, 2).Value = "12:28:16"
$wsTemp.Cells.Item(266, 3).Value = "12:00"
$wsTemp.Cells.Item(266, 4).Value = "Bathroom"
$wsTemp.Cells.Item(266, 5).Value = "22.9C"
$wsTemp.Cells.Item(266, 6).Value = "Active"
$wsTemp.Cells.Item(267, 1).NumberFormat = "@"
$wsTemp.Cells.Item(267, 1).Value = "2026-01-28"
$wsTemp.Cells.Item(267, 2).Value = "12:28:20"
$wsTemp.Cells.Item(267, 3).Value = "12:00"
$wsTemp.Cells.Item(267, 4).Value = "Bathroom"
$wsTemp.Cells.Item(267, 5).Value = "22.9C"
$wsTemp.Cells.Item(267, 6).Value = "Active"
$wsTemp.Cells.Item(268, 1).NumberFormat = "@"
$wsTemp.Cells.Item(268, 1).Value = "2026-01-28"
$wsTemp.Cells.Item(268, 2).Value = "12:28:24"
$wsTemp.Cells.Item(268, 3).Value = "12:00"
$wsTemp.Cells.Item(268, 4).Value = "Bathroom"
$wsTemp.Cells.Item(268, 5).Value = "22.9C"
$wsTemp.Cells.Item(268, 6).Value = "Active"
$wsTemp.Cells.Item(269, 1).NumberFormat = "@"
$wsTemp.Cells.Item(269, 1).Value = "2026-01-28"
$wsTemp.Cells.Item(269, 2).Value = "12:28:28"
$wsTemp.Cells.Item(269, 3).Value = "12:00"
$wsTemp.Cells.Item(269, 4).Value = "Bathroom"
$wsTemp.Cells.Item(269, 5).Value = "22.9C"
$wsTemp.Cells.Item(269, 6).Value = "Active"
$wsTemp.Cells.Item(270, 1).NumberFormat = "@"
$wsTemp.Cells.Item(270, 1).Value = "2026-01-28"
$wsTemp.Cells.Item(270, 2).Value = "12:28:36"
$wsTemp.Cells.Item(270, 3).Value = "12:00"
$wsTemp.Cells.Item(270, 4).Value = "Bathroom"
$wsTemp.Cells.Item(270, 5).Value = "22.9C"
$wsTemp.Cells.Item(270, 6).Value = "Active"
$wsTemp.Cells.Item(271, 1).NumberFormat = "@"
$wsTemp.Cells.Item(271, 1).Value = "2026-01-28"
$wsTemp.Cells.Item(271, 2).Value = "12:28:40"
$wsTemp.Cells.Item(271, 3).Value = "12:00"
$wsTemp.Cells.Item(271, 4).Value = "Bathroom"
$wsTemp.Cells.Item(271, 5).Value = "22.9C"
$wsTemp.Cells.Item(271, 6).Value = "Active"
$wsTemp.Cells.Item(272, 1).NumberFormat = "@"
$wsTemp.Cells.Item(272, 1).Value = "2026-01-28"
$wsTemp.Cells.Item(272, 2).Value = "12:28:44"
$wsTemp.Cells.Item(272, 3).Value = "12:00"
$wsTemp.Cells.Item(272, 4).Value = "Bathroom"
$wsTemp.Cells.Item(272, 5).Value = "22.9C"
$wsTemp.Cells.Item(272, 6).Value = "Active"
$wsTemp.Cells.Item(273, 1).NumberFormat = "@"
$wsTemp.Cells.Item(273, 1).Value = "2026-01-28"
$wsTemp.Cells.Item(273, 2).Value = "12:28:48"
$wsTemp.Cells.Item(273, 3).Value = "12:00"
$wsTemp.Cells.Item(273, 4).Value = "Bathroom"
$wsTemp.Cells.Item(273, 5).Value = "22.9C"
$wsTemp.Cells.Item(273, 6).Value = "Active"
$wsTemp.Cells.Item(274, 1).NumberFormat = "@"
$wsTemp.Cells.Item(274, 1).Value = "2026-01-28"
$wsTemp.Cells.Item(274, 2).Value = "12:28:52"
$wsTemp.Cells.Item(274, 3).Value = "12:00"
$wsTemp.Cells.Item(274, 4).Value = "Bathroom"
$wsTemp.Cells.Item(274, 5).Value = "22.9C"
$wsTemp.Cells.Item(274, 6).Value = "Active"
